$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 2-5 in columns A (WIID) and E (Date)
$ws.Range("A2").Value = 147455
$ws.Range("E2").Value = 43480

$ws.Range("A3").Value = 211405
$ws.Range("E3").Value = 43211

$ws.Range("A4").Value = 310925
$ws.Range("E4").Value = 42877

$ws.Range("A5").Value = 240995
$ws.Range("E5").Value = 43244

# Append new rows 6-9 with the same Description/Type/Status text as existing rows
$ws.Range("A6").Value = 477185
$ws.Range("B6").Value = "Calculate Client Security Hash"
$ws.Range("C6").Value = "WI5"
$ws.Range("D6").Value = "Open"
$ws.Range("E6").Value = 43433

$ws.Range("A7").Value = 496005
$ws.Range("B7").Value = "Calculate Client Security Hash"
$ws.Range("C7").Value = "WI5"
$ws.Range("D7").Value = "Open"
$ws.Range("E7").Value = 43495

$ws.Range("A8").Value = 393865
$ws.Range("B8").Value = "Calculate Client Security Hash"
$ws.Range("C8").Value = "WI5"
$ws.Range("D8").Value = "Open"
$ws.Range("E8").Value = 42775

$ws.Range("A9").Value = 217285
$ws.Range("B9").Value = "Calculate Client Security Hash"
$ws.Range("C9").Value = "WI5"
$ws.Range("D9").Value = "Open"
$ws.Range("E9").Value = 42866

# Copy the date number-format from an existing formatted date cell onto the
# newly added date cells so they share the same style (numFmt 14, m/d/yyyy)
# instead of creating brand new style entries.
$ws.Range("E2").Copy()
$ws.Range("E6:E9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the selection left behind on the sheet (whole rows 2:5 selected)
$ws.Range("A2:A5").EntireRow.Select()
